# Applies the update to the "AYKO" tracking sheet:
# two closed/removed cases are deleted from the table:
#   - old row 62 ("-406" / Olof palme 4144)
#   - old row 88 ("-536" / Olof palme 4142), which is the last row of the sheet
# Deleting them shifts the remaining rows up, so the data that used to be on
# rows 63..87 now lives on rows 62..86, and the sheet's last used row becomes 86.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AYKO")

# Delete the later row first so the earlier row number (62) still refers to
# the same original record when we delete it afterwards.
$ws.Rows.Item(88).Delete()
$ws.Rows.Item(62).Delete()
